# The "Coordenadas" sheet used to hold latitude/longitude pairs (as
# hyperlinked lookups to geodatos.net) for each province capital. The
# commit ("Lectura de datos de excel para mapa implementada") replaces
# those coordinates with plain integer values read for the map, and the
# now-stale hyperlinks are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coordenadas")

# New B/C values per row (2..25), replacing the old lat/long floats.
$newValues = @(
    @(354, 348),
    @(238, 266),
    @(349, 169),
    @(357, 131),
    @(365, 357),
    @(172, 216),
    @(128, 305),
    @(139, 455),
    @(312, 276),
    @(411, 166),
    @(210, 564),
    @(341, 163),
    @(117, 781),
    @(194, 191),
    @(203, 151),
    @(203, 81),
    @(203, 98),
    @(133, 271),
    @(178, 316),
    @(303, 270),
    @(228, 396),
    @(224, 175),
    @(139, 855),
    @(250, 504)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $i + 2
    $pair = $newValues[$i]
    $ws.Cells.Item($row, 2).Value = $pair[0]
    $ws.Cells.Item($row, 3).Value = $pair[1]
}

# The old values were rendered as hyperlinks (to display the number with a
# comma-decimal) - those no longer make sense for the new integer data.
$ws.Hyperlinks.Delete()

# Reflect where the user ended up after the edit.
$ws.Activate() | Out-Null
$ws.Range("H12").Select() | Out-Null
